# Rename the header row (row 1) to reflect the new "RSP" naming scheme
# and restore the plain (no leading-space) labels for the spectrum/talent
# columns. All other data rows keep referencing the same underlying text,
# so the shared-string table is rebuilt by the engine automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "RSP Pos. Ranking"
$ws.Range("H1").Value = "RSP Notes"
$ws.Range("E1").Value = "Comparison Spectrum"
$ws.Range("F1").Value = "Depth of Talent Score"
$ws.Range("G1").Value = "Depth of Talent Description"

# Update the window/selection state to match the saved view.
$ws.Range("G1").Select()
